# Remove the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph and the
# "© 2020 . Contact: ..." copyright/footer paragraph that followed the
# "LOQ4003: ..." requirement line, while keeping the surrounding blank
# paragraphs intact (matches the site rebuild diff that dropped the footer
# scraped content).

$d = $word.ActiveDocument

$start = $null
$end = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ($t -like "Ver no Jupiter*") {
        $start = $p.Range.Start
    }
    if ($t -like "*Powered by Jekyll and Github pages*") {
        $end = $p.Range.End
    }
}

if ($start -ne $null -and $end -ne $null) {
    $r = $d.Range($start, $end)
    $r.Delete()
}
